# draft_to_dos.xlsx - "still unidentified - gotta check procedures tomorra"
#
# 1) Mark two existing to-do items as done (strike-through the text).
# 2) Add a little "App" sub-section (two struck-through bullets) right
#    after the current list, then leave a blank row and append five new
#    open to-do items below it.
# 3) Leave the cursor on the last "App" bullet, like the author did.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- mark completed items ------------------------------------------------
# "Incorporate in the optimal Taylor rule section ..." -> done
$ws.Range("A3").Font.Strikethrough = $true
# "write out target criterion, expanding the terms so the signs become clear" -> done
$ws.Range("A6").Font.Strikethrough = $true

# --- new "open" to-dos appended at the bottom (rows 26-29, then 30) ------
# Written first so they claim shared-string slots 24-27 before the two
# "App" bullets below take 28-29, matching the order new items were typed.
$ws.Range("A26").Value = "get estimation identified"
$ws.Range("A27").Value = "redo PEA-VFI: figures, interpretations"
$ws.Range("A28").Value = "redo optimal Taylor rule: figures, table, interpretations"
$ws.Range("A29").Value = "redo IRFs of model in app, possibly interpretations"

# --- new "App" bullets (rows 23-24), already struck through ---------------
$ws.Range("A23").Value = "App alternative specifications of anchoring function"
$ws.Range("A23").Font.Strikethrough = $true

$ws.Range("A24").Value = "App estimation procedure"
$ws.Range("A24").Font.Strikethrough = $true

# last new open to-do, row 30
$ws.Range("A30").Value = "do a welfare bit that compares welfare under the optimal policy, an optimal TR under anchoring and an optimal TR under RE (may not make it into first draft)"

# --- leave selection where the author left it ------------------------------
[void]$ws.Range("A24").Select()
